# Dividend Calculation.xlsx - update Yearly "November" Taxable Account dividend
# and let dependent totals recalc; move the active tab/selection from Yearly
# to All Time as captured by the commit's saved view state.

$wb = $excel.ActiveWorkbook

$yearly = $wb.Worksheets.Item("Yearly")
$allTime = $wb.Worksheets.Item("All Time")

# --- Data edit -----------------------------------------------------------
# Yearly!L8 (2017 / November / Taxable Account) grew from 118.95 to 142.2.
# O8 (=SUM(L8:N8)), L15/O15 (=SUM(...)) on "Yearly" and the mirrored
# F8/I8/F46/I46 totals on "All Time" are formulas, so they recalc
# automatically once the source cell changes.
$yearly.Range("L8").Value = 142.2

# --- View / selection state ----------------------------------------------
# Originally "Yearly" was the active tab with K22 selected and "All Time"
# was scrolled to A25 with A55 selected. After the edit, "All Time" becomes
# the active tab (scrolled so row 19 is visible, N29 selected) and "Yearly"
# keeps a plain H21 selection.
$yearly.Range("H21").Select()

$allTime.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$allTime.Range("N29").Select()
